$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value2 = 42602.583657407406
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = "Random"

$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 69
$ws.Cells.Item($row, 9).Value = 31
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 58
$ws.Cells.Item($row, 13).Value = 42
